$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.301.57'
$ws.Range("E2").Value = '  -0.49%  '

$ws.Range("D3").Value = '2.604.36'
$ws.Range("E3").Value = '  +2.91%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.32'
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.42'
$ws.Range("E6").Value = '  -3.74%  '

$ws.Range("E7").Value = '  -1.13%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.576'
$ws.Range("E9").Value = '  +1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.27'
$ws.Range("E10").Value = '  +0.94%  '

$ws.Range("E11").Value = '  -0.63%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0839'
$ws.Range("E12").Value = '  +1.80%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.09'
$ws.Range("E13").Value = '  +2.40%  '

$ws.Range("D14").Value = '3.001.48'
$ws.Range("E14").Value = '  +2.63%  '

$ws.Range("E15").Value = '  +0.76%  '

$ws.Range("D16").Value = '2.610.54'
$ws.Range("E16").Value = '  +2.53%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.916'
$ws.Range("E17").Value = '  +2.55%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.89'
$ws.Range("E18").Value = '  -0.42%  '

$ws.Range("D19").Value = '46.366.64'
$ws.Range("E19").Value = '  -0.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000100'
$ws.Range("E20").Value = '  +1.83%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.90'
$ws.Range("E21").Value = '  -7.87%  '

$ws.Range("E22").Value = '  +1.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '71.20'
$ws.Range("E23").Value = '  +2.34%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '271.86'
$ws.Range("E24").Value = '  +7.08%  '

$ws.Range("E25").Value = '  +1.64%  '

$ws.Range("E26").Value = '  +2.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '29.34'
$ws.Range("E27").Value = '  +20.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  +0.10%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.01'
$ws.Range("E29").Value = '  -0.71%  '

$ws.Range("E30").Value = '  +1.12%  '

$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.21'
$ws.Range("E31").Value = '  -2.60%  '

$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '38.14'
$ws.Range("E32").Value = '  -8.88%  '

$ws.Range("E33").Value = '  +5.41%  '

$ws.Range("E34").Value = '  -3.01%  '

$ws.Range("E35").Value = '  -2.36%  '

$ws.Range("E36").Value = '  +2.04%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0832'
$ws.Range("E37").Value = '  -1.02%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '151.35'
$ws.Range("E38").Value = '  +1.00%  '

$ws.Range("E39").Value = '  +1.99%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.122'
$ws.Range("E40").Value = '  +1.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.19'
$ws.Range("E41").Value = '  +32.78%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '15.80'
$ws.Range("E42").Value = '  -3.54%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0328'
$ws.Range("E43").Value = '  +1.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.58'
$ws.Range("E44").Value = '  +1.27%  '

$ws.Range("E45").Value = '  -5.29%  '

$ws.Range("D46").Value = '2.116.00'
$ws.Range("E46").Value = '  +6.05%  '

$ws.Range("E47").Value = '  -0.12%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '93.48'
$ws.Range("E48").Value = '  -0.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.51'
$ws.Range("E49").Value = '  +7.43%  '

$ws.Range("E50").Value = '  -5.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '108.24'
$ws.Range("E51").Value = '  +1.81%  '
